$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 525, shifting existing rows 525:653 down to 526:654
$ws.Rows.Item(525).Insert()

# Populate the newly inserted row 525 with the new record's data
$ws.Cells.Item(525, 1).Value = 10
$ws.Cells.Item(525, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(525, 3).Value = "La Araucanía"
$ws.Cells.Item(525, 4).Value = 45204
$ws.Cells.Item(525, 5).Value = 9
$ws.Cells.Item(525, 6).Value = 100112040
$ws.Cells.Item(525, 7).Value = "Cilantro"
$ws.Cells.Item(525, 8).Value = "Sin especificar"
$ws.Cells.Item(525, 9).Value = "Primera"
$ws.Cells.Item(525, 10).Value = 125
$ws.Cells.Item(525, 11).Value = 3600
$ws.Cells.Item(525, 12).Value = 3600
$ws.Cells.Item(525, 13).Value = 3600
$ws.Cells.Item(525, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(525, 15).Value = "Región Metropolitana"
$ws.Cells.Item(525, 16).Value = 1800
$ws.Cells.Item(525, 17).Value = 2
$ws.Cells.Item(525, 18).Value = "Hortaliza"

# Make sure the new row uses the same date-number-format style as the other rows in column D
$ws.Cells.Item(525, 4).NumberFormat = $ws.Cells.Item(526, 4).NumberFormat
